# Add the 8th function (linear_func) block: rows 32-35, mirroring the
# structure of the existing function blocks (e.g. rows 28-31 / poly_func).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the previous block's (poly_func, rows 28-31) cell formatting + merge
# down onto the new rows 32-35, then re-apply the uniform thin border so the
# "Estimator" label column keeps the same bold/centered/bordered look as
# every other block (A4:A7, ..., A28:A31) instead of the visually-merged
# per-edge border split that a straight copy of a merged range leaves behind.
$ws.Range("A28:A31").Copy($ws.Range("A32:A35"))
$ws.Range("B28:B31").Copy($ws.Range("B32:B35"))
$ws.Range("A32:A35").Borders.LineStyle = 1

# --- Row 32 ---
$ws.Range("A32").Value = "linear_func"
$ws.Range("B32").Value = "co1"
$ws.Range("C32").Value2 = -0.7090672712095367
$ws.Range("D32").Value2 = 0.7051408412111475
$ws.Range("E32").Value2 = 0.0005270846108141711
$ws.Range("F32").Value2 = -0.07585097507128481
$ws.Range("G32").Value2 = -0.02323376716424582
$ws.Range("H32").Value2 = 1.416035604360695
$ws.Range("J32").Value2 = -1.080528565795282
$ws.Range("K32").Value2 = 1.074592175897344
$ws.Range("L32").Value2 = 0.0005642027868387669
$ws.Range("M32").Value2 = -0.07550442330232486
$ws.Range("N32").Value2 = -0.02307341806374083
$ws.Range("O32").Value2 = 0.9289311343387469

# --- Row 33 ---
$ws.Range("B33").Value = "co2"
$ws.Range("C33").Value2 = -0.7516685886265817
$ws.Range("D33").Value2 = 0.6595410402246661
$ws.Range("E33").Value2 = 0.06671850804371433
$ws.Range("F33").Value2 = -0.2487130079975295
$ws.Range("G33").Value2 = 0.003689032806593822
$ws.Range("H33").Value2 = 1.047740175148692
$ws.Range("J33").Value2 = -0.8574933893933496
$ws.Range("K33").Value2 = 0.7524733899579436
$ws.Range("L33").Value2 = 0.06674306684259475
$ws.Range("M33").Value2 = -0.2491077528920168
$ws.Range("N33").Value2 = 0.003692653977174958
$ws.Range("O33").Value2 = 0.919095041468339

# --- Row 34 ---
$ws.Range("B34").Value = "co3"
$ws.Range("C34").Value2 = -0.9999938087821889
$ws.Range("D34").Value2 = -0.003523743343453054
$ws.Range("E34").Value2 = 0.05901272340870947
$ws.Range("F34").Value2 = -0.2239929185852499
$ws.Range("G34").Value2 = -0.08319682079167365
$ws.Range("H34").Value2 = 0.02634630246596919
$ws.Range("J34").Value2 = -0.7800501824005512
$ws.Range("K34").Value2 = -0.003461882162105671
$ws.Range("L34").Value2 = 0.05901831382857809
$ws.Range("M34").Value2 = -0.2239422745182012
$ws.Range("N34").Value2 = -0.08323356207156772
$ws.Range("O34").Value2 = 0.03376579041968644

# --- Row 35 ---
$ws.Range("B35").Value = "co4"
$ws.Range("C35").Value2 = 0.4766520604249416
$ws.Range("D35").Value2 = -0.8790920648163564
$ws.Range("E35").Value2 = 0.07648983205004396
$ws.Range("F35").Value2 = -0.09020498935921953
$ws.Range("G35").Value2 = 0.02352424865596539
$ws.Range("H35").Value2 = 0.5277006196296115
$ws.Range("J35").Value2 = -1.472723895481876
$ws.Range("K35").Value2 = -0.02069103131645775
$ws.Range("L35").Value2 = 0.07287479566808779
$ws.Range("M35").Value2 = -0.1111105239325738
$ws.Range("N35").Value2 = 0.02486304463173006
$ws.Range("O35").Value2 = 0.1137925998256436
